$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; this shifts existing rows 31..61 down to 32..62
# and extends the used range to A1:R62, matching the target diff.
$ws.Rows(31).Insert()

# Populate the newly inserted row 31 with the new weekly price record.
$ws.Cells.Item(31, 1).Value = 10
$ws.Cells.Item(31, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(31, 3).Value = "La Araucanía"
$ws.Cells.Item(31, 4).Value = 44638
$ws.Cells.Item(31, 5).Value = 9
$ws.Cells.Item(31, 6).Value = 100114002
$ws.Cells.Item(31, 7).Value = "Camote"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 40
$ws.Cells.Item(31, 11).Value = 18000
$ws.Cells.Item(31, 12).Value = 18000
$ws.Cells.Item(31, 13).Value = 18000
$ws.Cells.Item(31, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(31, 15).Value = "Perú"
$ws.Cells.Item(31, 16).Value = 900
$ws.Cells.Item(31, 17).Value = 20
$ws.Cells.Item(31, 18).Value = "Hortaliza"
